$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the unit_rate column (D) values from rupee to lakh by dividing by
# the appropriate factor, expressed as formulas so Excel recalculates them.
$ws.Range("D2").Formula = "=5100/10000"
$ws.Range("D3").Formula = "=3600/100000"
$ws.Range("D4").Formula = "=550/100000"
$ws.Range("D5").Formula = "=380/100000"
$ws.Range("D6").Formula = "=7.8/100000"

# Widen column D to fit the new note text.
$ws.Columns.Item(4).ColumnWidth = 25.83

# Add a footnote explaining the unit conversion.
$ws.Range("D7").Value = "in lakh (converted from rupee)"

# Update the active selection to reflect where the user left off editing.
$ws.Range("D7").Select()
